$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.396.41'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '1.722.80'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.97'
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9987'
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4887'
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2600'
$ws.Range("E8").Value = '  -2.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06183'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = '1.726.06'
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06979'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.52'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.513'
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5980'
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.00'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9989'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '26.388.58'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9984'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007136'
$ws.Range("E19").Value = '  -1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.29'
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").Value = '1.948.03'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.442'
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.496'
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.093'
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.09'
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.23'
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.396'
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '106.25'
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.733'
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.903'
$ws.Range("E30").Value = '  -1.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08032'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04491'
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9956'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  -2.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9176'
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.960'
$ws.Range("E38").Value = '  -3.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.380'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9982'
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("E41").Value = '  -2.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.89'
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.440'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3839'
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.893'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1161'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05362'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.26'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.656'
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.06'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.217'
$ws.Range("E51").Value = '  -2.74%  '
